# "Actualización desde MV -datos-"
#
# The source adds one more reporting period to the PIB (GDP) series table:
# a new last column "BH" labeled "Agosto.2021", placed right after the
# current last column "BG" ("Mayo.2021"). Like every previous "new period"
# column before it, it starts out just restating each row's latest known
# figure (the BG value), matching the source file's pattern where a newly
# opened period column is a copy of the prior one until revised data comes
# in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: clone BG1's formatting (bold, centered, thin border) onto the new
# BH1 cell, then give it its own caption.
$ws.Range("BG1").Copy()
$ws.Range("BH1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("BH1").Value = "Agosto.2021"

# Data rows 2-19: BH repeats BG's value for every series.
$ws.Range("BH2:BH19").Value2 = $ws.Range("BG2:BG19").Value2
